$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> FAPs) ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.183046666666667
$ws.Range("N2").Value = 3.54914
$ws.Range("O2").Value = 0.6222589862820888
$ws.Range("P2").Value = 0.6222589862820888
$ws.Range("Q2").Value = 7.382338969039998
$ws.Range("R2").Value = 66.44105072135999
$ws.Range("S2").Value = 0.01077896009584504
$ws.Range("T2").Value = 0.01077896009584504

# --- Row 3 (ECs -> MuSCs) ---
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("M3").Value = 0.718166
$ws.Range("N3").Value = 2.154498
$ws.Range("O3").Value = 0.3777410137179113
$ws.Range("P3").Value = 0.3777410137179112
$ws.Range("Q3").Value = 4.481433401927999
$ws.Range("R3").Value = 40.33290061735199
$ws.Range("S3").Value = 0.006543345139548721
$ws.Range("T3").Value = 0.006543345139548721

# --- Row 4 (FAPs -> FAPs) ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 345.566579
$ws.Range("H4").Value = 1036.699737
$ws.Range("I4").Value = 0.9592798330716089
$ws.Range("J4").Value = 0.9592798330716091
$ws.Range("M4").Value = 1.183046666666667
$ws.Range("N4").Value = 3.54914
$ws.Range("O4").Value = 0.6222589862820888
$ws.Range("P4").Value = 0.6222589862820888
$ws.Range("Q4").Value = 408.8213893973533
$ws.Range("R4").Value = 3679.39250457618
$ws.Range("S4").Value = 0.5969204964879907
$ws.Range("T4").Value = 0.5969204964879908

# --- Row 5 (FAPs -> MuSCs) ---
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.718166
$ws.Range("N5").Value = 2.154498
$ws.Range("O5").Value = 0.3777410137179113
$ws.Range("P5").Value = 0.3777410137179112
$ws.Range("Q5").Value = 248.174167774114
$ws.Range("R5").Value = 2233.567509967026
$ws.Range("S5").Value = 0.3623593365836183
$ws.Range("T5").Value = 0.3623593365836183

# --- Row 6 (MuSCs -> FAPs) ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("G6").Value = 8.428738666666666
$ws.Range("H6").Value = 25.286216
$ws.Range("I6").Value = 0.02339786169299727
$ws.Range("J6").Value = 0.02339786169299728
$ws.Range("O6").Value = 0.6222589862820888
$ws.Range("P6").Value = 0.6222589862820888
$ws.Range("Q6").Value = 9.971591183804442
$ws.Range("R6").Value = 89.74432065424
$ws.Range("S6").Value = 0.014559529698253
$ws.Range("T6").Value = 0.014559529698253

# --- Row 7 (MuSCs -> MuSCs) ---
$ws.Range("A7").Value = "MuSCs"
$ws.Range("G7").Value = 8.428738666666666
$ws.Range("H7").Value = 25.286216
$ws.Range("I7").Value = 0.02339786169299727
$ws.Range("J7").Value = 0.02339786169299728
$ws.Range("O7").Value = 0.3777410137179113
$ws.Range("P7").Value = 0.3777410137179112
$ws.Range("Q7").Value = 6.053233533285333
$ws.Range("R7").Value = 54.47910179956799
$ws.Range("S7").Value = 0.008838331994744273
$ws.Range("T7").Value = 0.008838331994744274

# Rows 8,9,10 (MuSCs -> * ) are removed entirely, the table now only has 6 data rows
$ws.Rows("8:10").Delete()
